$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Update the "Port Config (tab name)" values and the dependent counts
# Row 3 (Run ID 2)
$ws.Range("B3").Value = 20
$ws.Range("K3").Value = "config2"

# Row 4 (Run ID 3)
$ws.Range("B4").Value = 30
$ws.Range("K4").Value = "config3"

# Row 5 (Run ID 4)
$ws.Range("B5").Value = 40
$ws.Range("K5").Value = "config4"

# Row 6 (Run ID 5)
$ws.Range("B6").Value = 50
$ws.Range("K6").Value = "config5"

# Row 8 (Run ID 7)
$ws.Range("C8").Value = 20
$ws.Range("K8").Value = "config2"

# Row 9 (Run ID 8)
$ws.Range("C9").Value = 30
$ws.Range("K9").Value = "config3"

# Row 10 (Run ID 9)
$ws.Range("C10").Value = 40
$ws.Range("K10").Value = "config4"

# Row 11 (Run ID 10)
$ws.Range("C11").Value = 50
$ws.Range("K11").Value = "config5"

# Row 12 (Run ID 11)
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 0

# Row 13 (Run ID 12)
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = 10
$ws.Range("K13").Value = "config2"

# Row 14 (Run ID 13)
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 15
$ws.Range("K14").Value = "config3"

# Row 15 (Run ID 14)
$ws.Range("B15").Value = 20
$ws.Range("C15").Value = 20
$ws.Range("K15").Value = "config4"

# Row 16 (Run ID 15)
$ws.Range("B16").Value = 25
$ws.Range("C16").Value = 25
$ws.Range("K16").Value = "config5"

# Make the "params" sheet the active tab, with K1:K16 selected
$ws.Activate() | Out-Null
$ws.Range("K1:K16").Select() | Out-Null
